$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest scraped values.
# Cells whose new value would otherwise be auto-detected as a number by Excel
# are first switched to Text format so they keep the exact original string look
# (leading/trailing zeros, thousands-dot grouping, etc.), matching column D's
# existing all-text convention.

$ws.Range("D2").Value = "40.099.26"
$ws.Range("E2").Value = "  +2.90%  "
$ws.Range("D3").Value = "2.238.92"
$ws.Range("E3").Value = "  +0.80%  "
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "295.59"
$ws.Range("E5").Value = "  -0.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "86.95"
$ws.Range("E6").Value = "  +8.48%  "
$ws.Range("E7").Value = "  +2.46%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.474"
$ws.Range("E9").Value = "  +3.53%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "31.40"
$ws.Range("E10").Value = "  +12.06%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0795"
$ws.Range("E11").Value = "  +2.68%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.21"
$ws.Range("E12").Value = "  +3.21%  "
$ws.Range("E13").Value = "  +0.81%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.50"
$ws.Range("E14").Value = "  +6.38%  "
$ws.Range("D15").Value = "2.590.82"
$ws.Range("E15").Value = "  +1.35%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.19"
$ws.Range("E16").Value = "  +1.87%  "
$ws.Range("D17").Value = "2.237.15"
$ws.Range("E17").Value = "  -0.44%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.731"
$ws.Range("E18").Value = "  +2.69%  "
$ws.Range("D19").Value = "40.045.99"
$ws.Range("E19").Value = "  +3.06%  "
$ws.Range("D20").Value = "0.0₃0892"
$ws.Range("E20").Value = "  +3.96%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.82"
$ws.Range("E21").Value = "  +1.65%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.97"
$ws.Range("E22").Value = "  +11.85%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.62"
$ws.Range("E23").Value = "  +1.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "235.78"
$ws.Range("E24").Value = "  +4.55%  "
$ws.Range("E25").Value = "  -0.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.47"
$ws.Range("E26").Value = "  +3.90%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.85"
$ws.Range("E27").Value = "  +5.72%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.88"
$ws.Range("E28").Value = "  +2.64%  "
$ws.Range("E29").Value = "  +5.43%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.24"
$ws.Range("E30").Value = "  +3.42%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "33.51"
$ws.Range("E31").Value = "  +7.66%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "152.53"
$ws.Range("E32").Value = "  +3.27%  "
$ws.Range("E33").Value = "  +0.16%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.89"
$ws.Range("E34").Value = "  +3.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0720"
$ws.Range("E35").Value = "  +4.68%  "
$ws.Range("E36").Value = "  +3.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "16.52"
$ws.Range("E37").Value = "  +13.74%  "
$ws.Range("E38").Value = "  +3.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.100"
$ws.Range("E39").Value = "  +5.71%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.72"
$ws.Range("E40").Value = "  +2.40%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.70"
$ws.Range("E41").Value = "  +7.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.86"
$ws.Range("E42").Value = "  +7.58%  "
$ws.Range("D43").Value = "2.038.94"
$ws.Range("E43").Value = "  +6.82%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.27"
$ws.Range("E44").Value = "  +9.39%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0272"
$ws.Range("E45").Value = "  +7.16%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.03"
$ws.Range("E46").Value = "  +13.15%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "16.30"
$ws.Range("E47").Value = "  +0.95%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.58"
$ws.Range("E48").Value = "  +2.54%  "
$ws.Range("D49").Value = "2.460.36"
$ws.Range("E49").Value = "  +1.59%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "70.94"
$ws.Range("E50").Value = "  +0.27%  "
$ws.Range("E51").Value = "  +14.77%  "
